$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Jengibre (Terminal La Palmera de La Serena)
# is inserted as the first data row (row 148), pushing the existing
# records (old rows 148-165) down by one (new rows 149-166).
$ws.Rows.Item(148).Insert()

$ws.Cells.Item(148, 1).Value = 8
$ws.Cells.Item(148, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(148, 3).Value = "Coquimbo"
$ws.Cells.Item(148, 4).Value = 45142
$ws.Cells.Item(148, 5).Value = 4
$ws.Cells.Item(148, 6).Value = 100114007
$ws.Cells.Item(148, 7).Value = "Jengibre"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 360
$ws.Cells.Item(148, 11).Value = 18000
$ws.Cells.Item(148, 12).Value = 19000
$ws.Cells.Item(148, 13).Value = 18500
$ws.Cells.Item(148, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(148, 15).Value = "Perú"
$ws.Cells.Item(148, 16).Value = 1423
$ws.Cells.Item(148, 17).Value = 13
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Match the date-format style used by the other rows in column D.
$ws.Cells.Item(148, 4).NumberFormat = $ws.Cells.Item(149, 4).NumberFormat
